$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 107 (idx 0)
$ws.Range("H107").Value = 1249.5
$ws.Range("I107").Value = 973.6316
$ws.Range("K107").Value = 973.6316
$ws.Range("M107").Value = 946.3684
# Row 113 (idx 1)
$ws.Range("H113").Value = 5247.143
$ws.Range("I113").Value = 3152
$ws.Range("J113").Value = 6411.1113
$ws.Range("K113").Value = 3152
$ws.Range("L113").Value = 6411.1113
$ws.Range("M113").Value = 102
$ws.Range("N113").Value = -12919.1113
# Row 116 (idx 2)
$ws.Range("H116").Value = 3874.125
$ws.Range("I116").Value = 4200.2
$ws.Range("J116").Value = 3330.6667
$ws.Range("K116").Value = 4200.2
$ws.Range("L116").Value = 3330.6667
$ws.Range("M116").Value = -758.1999999999998
$ws.Range("N116").Value = -10214.6667
# Row 129 (idx 3)
$ws.Range("H129").Value = 1066.5
$ws.Range("J129").Value = 1161.1111
$ws.Range("L129").Value = 3483.3333
$ws.Range("N129").Value = -13483.3333
# Row 132 (idx 4)
$ws.Range("H132").Value = 3111.4062
$ws.Range("I132").Value = 1802.25
$ws.Range("J132").Value = 7038.875
$ws.Range("K132").Value = 5406.75
$ws.Range("L132").Value = 21116.625
$ws.Range("M132").Value = -2876.75
$ws.Range("N132").Value = -26176.625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (idx 5)
$ws.Range("H45").Value = 1228.2
$ws.Range("I45").Value = 1012.3
$ws.Range("J45").Value = 1660
$ws.Range("K45").Value = 1012.3
$ws.Range("L45").Value = 1660
$ws.Range("M45").Value = -635.3
$ws.Range("N45").Value = -2414
# Row 61 (idx 6)
$ws.Range("H61").Value = 2749.8684
$ws.Range("I61").Value = 1541.8572
$ws.Range("K61").Value = 1541.8572
$ws.Range("M61").Value = -1329.8572
# Row 132 (idx 7)
$ws.Range("H132").Value = 50064.316
$ws.Range("I132").Value = 96733.63
$ws.Range("J132").Value = 3395
$ws.Range("K132").Value = 290200.89
$ws.Range("L132").Value = 10185
$ws.Range("M132").Value = -287670.89
$ws.Range("N132").Value = -15245
# Row 136 (idx 8)
$ws.Range("H136").Value = 2749.8684
$ws.Range("I136").Value = 1541.8572
$ws.Range("K136").Value = 4625.571599999999
$ws.Range("M136").Value = -2075.571599999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (idx 9)
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1574
# Row 31 (idx 10)
$ws.Range("H31").Value = 2100.6428
$ws.Range("I31").Value = 1141.1364
$ws.Range("J31").Value = 5618.8335
$ws.Range("K31").Value = 1141.1364
$ws.Range("L31").Value = 5618.8335
$ws.Range("M31").Value = -846.1364000000001
$ws.Range("N31").Value = -6208.8335
# Row 34 (idx 11)
$ws.Range("H34").Value = 2100.6428
$ws.Range("I34").Value = 1141.1364
$ws.Range("J34").Value = 5618.8335
$ws.Range("K34").Value = 1141.1364
$ws.Range("L34").Value = 5618.8335
$ws.Range("M34").Value = -939.1364000000001
$ws.Range("N34").Value = -6022.8335
# Row 99 (idx 12)
$ws.Range("H99").Value = 78988.30499999999
$ws.Range("I99").Value = 85200.336
$ws.Range("K99").Value = 85200.336
$ws.Range("M99").Value = -83702.336
# Row 113 (idx 13)
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5340
# Row 122 (idx 14)
$ws.Range("H122").Value = 2170.3157
$ws.Range("I122").Value = 2527.0715
$ws.Range("J122").Value = 1171.4
$ws.Range("K122").Value = 7581.2145
$ws.Range("L122").Value = 3514.2
$ws.Range("M122").Value = -5131.2145
$ws.Range("N122").Value = -8414.200000000001
# Row 126 (idx 15)
$ws.Range("H126").Value = 78988.30499999999
$ws.Range("I126").Value = 85200.336
$ws.Range("K126").Value = 255601.008
$ws.Range("M126").Value = -253131.008
# Row 132 (idx 16)
$ws.Range("H132").Value = 2458.3547
$ws.Range("I132").Value = 1539.3846
$ws.Range("J132").Value = 3122.0557
$ws.Range("K132").Value = 4618.1538
$ws.Range("L132").Value = 9366.167099999999
$ws.Range("M132").Value = -2088.1538
$ws.Range("N132").Value = -14426.1671

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12 (idx 17)
$ws.Range("H12").Value = 42.733334
$ws.Range("I12").Value = 27.25
$ws.Range("K12").Value = 81.75
$ws.Range("M12").Value = 91.25
# Row 97 (idx 18)
$ws.Range("H97").Value = 300.0909
$ws.Range("I97").Value = 150.8
$ws.Range("J97").Value = 424.5
$ws.Range("K97").Value = 452.4
$ws.Range("L97").Value = 1273.5
$ws.Range("M97").Value = 43.59999999999997
$ws.Range("N97").Value = -2265.5
# Row 113 (idx 19)
$ws.Range("H113").Value = 624.24
$ws.Range("I113").Value = 624.6875
$ws.Range("J113").Value = 623.44446
$ws.Range("K113").Value = 1874.0625
$ws.Range("L113").Value = 1870.33338
$ws.Range("M113").Value = 295.9375
$ws.Range("N113").Value = -6210.33338
# Row 122 (idx 20)
$ws.Range("H122").Value = 518.92
$ws.Range("I122").Value = 413.2857
$ws.Range("J122").Value = 1073.5
$ws.Range("K122").Value = 3719.5713
$ws.Range("L122").Value = 9661.5
$ws.Range("M122").Value = -1269.5713
$ws.Range("N122").Value = -14561.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (idx 21)
$ws.Range("H80").Value = 3186
$ws.Range("I80").Value = 3405.319
$ws.Range("J80").Value = 2670.6
$ws.Range("K80").Value = 3405.319
$ws.Range("L80").Value = 2670.6
$ws.Range("M80").Value = -2407.319
$ws.Range("N80").Value = -4666.6
# Row 83 (idx 22)
$ws.Range("H83").Value = 3186
$ws.Range("I83").Value = 3405.319
$ws.Range("J83").Value = 2670.6
$ws.Range("K83").Value = 17026.595
$ws.Range("L83").Value = 13353
$ws.Range("M83").Value = -12034.595
$ws.Range("N83").Value = -23337
# Row 126 (idx 23)
$ws.Range("H126").Value = 2858.1333
$ws.Range("J126").Value = 2998.5715
$ws.Range("L126").Value = 8995.7145
$ws.Range("N126").Value = -13935.7145

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (idx 24)
$ws.Range("H22").Value = 772.6667
$ws.Range("I22").Value = 613.3333
$ws.Range("J22").Value = 812.5
$ws.Range("K22").Value = 613.3333
$ws.Range("L22").Value = 812.5
$ws.Range("M22").Value = -318.3333
$ws.Range("N22").Value = -1402.5
# Row 27 (idx 25)
$ws.Range("H27").Value = 772.6667
$ws.Range("I27").Value = 613.3333
$ws.Range("J27").Value = 812.5
$ws.Range("K27").Value = 613.3333
$ws.Range("L27").Value = 812.5
$ws.Range("M27").Value = -506.3333
$ws.Range("N27").Value = -1026.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62 (idx 26)
$ws.Range("H62").Value = 4272.222
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4316.129
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4316.129
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5564.129
# Row 65 (idx 27)
$ws.Range("H65").Value = 4272.222
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4316.129
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 21580.645
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -27820.645
